$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Growth Temperature" column (F)
# so that it becomes "Shaking speed[rpm]" and the old column shifts to G.
$ws.Range("F1").EntireColumn.Insert()

# Insert another new column before the (now shifted) "Plate" column (H)
# so that it becomes "Humidity[%]" and the old "Plate" column shifts to I.
$ws.Range("H1").EntireColumn.Insert()

# Header row
$ws.Range("F1").Value = "Shaking speed[rpm]"
$ws.Range("H1").Value = "Humidity[%]"

# Fill the new data columns for every data row (2 through 49)
for ($row = 2; $row -le 49; $row++) {
    $ws.Cells.Item($row, 6).Value = 800
    $ws.Cells.Item($row, 8).Value = 90
}
